# Progress log update: append a new "Week 10" block (label row, three daily
# entries, and a totals row) after the existing "Week 9" block that ends at
# row 47.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bring over formatting for the new rows from existing, equivalent rows -
# Row 48 (new week label) gets its look from row 40 (the "Week 9" label row).
# Only column A is populated on a label row, so only copy that cell.
$ws.Range("A40").Copy()
$ws.Range("A48").PasteSpecial(-4122)

# Rows 49:51 (three daily entries) get their look from rows 41:43, which are
# plain (non-total) daily rows belonging to the "Week 9" block.
$ws.Range("A41:E43").Copy()
$ws.Range("A49:E51").PasteSpecial(-4122)

# Row 52 (week total) gets its look from row 47, the "Week 9" total row.
$ws.Range("A47:E47").Copy()
$ws.Range("A52:E52").PasteSpecial(-4122)

# The PasteSpecial calls above also copied formula/text contents along with
# the formatting; clear those back out so the cells start empty before we
# fill in the real values/formulas below (this also avoids stale cached
# calculation results for the new formulas).
$ws.Range("A48:E52").ClearContents()

# --- Week 10 label -----------------------------------------------------
$ws.Range("A48").Value = "Week 10"

# --- Daily entries -------------------------------------------------------
$ws.Range("A49").Value = 42280
$ws.Range("B49").Value = 0.875
$ws.Range("C49").Value = 0.10416666666666667
$ws.Range("D49").Value = 0
$ws.Range("E49").Formula = "=MOD(C49-B49,1)*24-D49"

$ws.Range("A50").Value = 42281
$ws.Range("B50").Value = 0.5
$ws.Range("C50").Value = 0.91666666666666663
$ws.Range("D50").Value = 6
$ws.Range("E50").Formula = "=MOD(C50-B50,1)*24-D50"

$ws.Range("A51").Value = 42282
$ws.Range("B51").Value = 0.47916666666666669
$ws.Range("C51").Value = 0.083333333333333329
$ws.Range("D51").Value = 1
$ws.Range("E51").Formula = "=MOD(C51-B51,1)*24-D51"

# --- Week total ------------------------------------------------------------
$ws.Range("D52").Value = "Total"
$ws.Range("E52").Formula = "=SUM(E49:E51)"

# --- Update the view so it shows the newly added rows, like the author's --
$win = $excel.ActiveWindow
$win.ScrollRow = 31
$win.ScrollColumn = 1
$ws.Range("C54").Select()
